$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "Wrong_Entity_NonEvent_as_Event"
$ws.Range("B2").Value = 65

# Update row 3
$ws.Range("A3").Value = "Wrong_Entity_Event_as_NonEvent"
$ws.Range("B3").Value = 28

# Update row 4
$ws.Range("B4").Value = 27

# Update row 5
$ws.Range("A5").Value = "Wrong_Tag_E_as_I"
$ws.Range("B5").Value = 3

# Add row 6, copying style from row 5
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)
$ws.Range("A6").Value = "Wrong_Tag_B_as_I"
$ws.Range("B6").Value = 1
